$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.338.60"
$ws.Range("E2").Value = "  -2.65%  "
$ws.Range("D3").Value = "3.682.95"
$ws.Range("E3").Value = "  -3.25%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "687.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.71%  "
$ws.Range("D7").Value = "3.679.94"
$ws.Range("E7").Value = "  -3.30%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.497"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.43%  "
$ws.Range("E10").Value = "  -8.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.36"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.440"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -7.24%  "
$ws.Range("E13").Value = "  -5.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -8.00%  "
$ws.Range("D15").Value = "4.306.76"
$ws.Range("E15").Value = "  -3.22%  "
$ws.Range("D16").Value = "3.684.89"
$ws.Range("E16").Value = "  -3.58%  "
$ws.Range("D17").Value = "69.410.85"
$ws.Range("E17").Value = "  -2.61%  "
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -9.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "476.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.659"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.98%  "
$ws.Range("D25").Value = "3.832.67"
$ws.Range("E25").Value = "  -3.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000128"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -9.77%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -9.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.80"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -10.91%  "
$ws.Range("E31").Value = "  -10.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.78"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.05"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.166"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.19%  "
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.84"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.65%  "
$ws.Range("D37").Value = "3.658.54"
$ws.Range("E37").Value = "  -2.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.37"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.28"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0918"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.99%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.951"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "163.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "48.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "29.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.64%  "
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.74"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -16.15%  "
$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000280"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -9.20%  "
$ws.Range("B51").Value = "ONDO"
$ws.Range("C51").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.45%  "
